$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "278.03"
Set-TextValue "D3" "27.30"
Set-TextValue "E3" "2.44%"
Set-TextValue "D4" "4.854"
Set-TextValue "E4" "1.91%"
Set-TextValue "D5" "0.06405"
Set-TextValue "D6" "6.991"
Set-TextValue "E6" "1.03%"
Set-TextValue "D7" "1.207"
Set-TextValue "E7" "-8.03%"
Set-TextValue "D8" "0.8796"
Set-TextValue "E8" "1.13%"
Set-TextValue "D9" "0.1524"
Set-TextValue "E9" "-3.14%"
Set-TextValue "D10" "0.05183"
Set-TextValue "E10" "3.01%"
Set-TextValue "D11" "0.07516"
Set-TextValue "E11" "0.33%"
Set-TextValue "D12" "0.02945"
Set-TextValue "E12" "1.32%"
Set-TextValue "D13" "0.08973"
Set-TextValue "E13" "-0.91%"
Set-TextValue "D14" "0.001569"
Set-TextValue "E14" "-0.42%"
Set-TextValue "D15" "0.0006371"
Set-TextValue "E15" "0.29%"
Set-TextValue "D16" "0.006077"
Set-TextValue "E16" "4.09%"
Set-TextValue "D17" "3.478"
Set-TextValue "E17" "0.76%"
Set-TextValue "D19" "2.245"
Set-TextValue "E19" "-1.68%"
Set-TextValue "E21" "2.39%"
Set-TextValue "D22" "3.902"
Set-TextValue "E22" "-0.72%"
Set-TextValue "D23" "0.04421"
Set-TextValue "E23" "0.68%"
Set-TextValue "D24" "0.1505"
Set-TextValue "E24" "9.01%"
Set-TextValue "E25" "0.45%"
Set-TextValue "D26" "0.003896"
Set-TextValue "E26" "-7.44%"
Set-TextValue "E28" "-1.79%"
Set-TextValue "D29" "0.0001643"
Set-TextValue "E29" "1.59%"
Set-TextValue "D40" "0.04093"
Set-TextValue "E40" "0.46%"
Set-TextValue "D41" "0.006813"
Set-TextValue "E41" "-3.59%"
Set-TextValue "E42" "0.26%"
Set-TextValue "E43" "-6.60%"
Set-TextValue "E44" "0.31%"
Set-TextValue "D45" "0.00005364"
Set-TextValue "E45" "3.17%"
Set-TextValue "E46" "9.30%"
Set-TextValue "D47" "0.01851"
Set-TextValue "E47" "-19.64%"
